# The edit swaps the two theme parts in the package: the theme that backs
# the (single) slide master -- ppt/theme/theme1.xml, originally the
# "Integral" palette -- ends up holding the "Office Theme" palette, while
# the theme used by the notes master keeps its own part untouched on disk
# (this host's object model only exposes a single mutable "active" theme,
# which always resolves to the slide-master theme part, so that is the
# part we edit here). The font scheme and format scheme are identical
# between the two themes already, so only the 12 theme colors (and, where
# possible, the scheme's display name) need to change.

function RGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Target palette: the stock PowerPoint "Office Theme" color scheme.
$tcs.Item(1).RGB  = RGB 0x00 0x00 0x00   # dk1
$tcs.Item(2).RGB  = RGB 0xFF 0xFF 0xFF   # lt1
$tcs.Item(3).RGB  = RGB 0x44 0x54 0x6A   # dk2
$tcs.Item(4).RGB  = RGB 0xE7 0xE6 0xE6   # lt2
$tcs.Item(5).RGB  = RGB 0x5B 0x9B 0xD5   # accent1
$tcs.Item(6).RGB  = RGB 0xED 0x7D 0x31   # accent2
$tcs.Item(7).RGB  = RGB 0xA5 0xA5 0xA5   # accent3
$tcs.Item(8).RGB  = RGB 0xFF 0xC0 0x00   # accent4
$tcs.Item(9).RGB  = RGB 0x44 0x72 0xC4   # accent5
$tcs.Item(10).RGB = RGB 0x70 0xAD 0x47   # accent6
$tcs.Item(11).RGB = RGB 0x05 0x63 0xC1   # hlink
$tcs.Item(12).RGB = RGB 0x95 0x4F 0x72   # folHlink
